$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the login records (token + id columns) for the three tracked
# users with newly minted UUIDs / JWTs, mirroring a fresh xlUtility login
# cycle (daniel5f, Jorge2525, mario35).

$ws.Range("D2").Value = "59b0c1bd-bedd-4aab-9078-d23b4fc02fee"
$ws.Range("D3").Value = "26a12013-bdef-475a-969f-7b5e053dbc02"
$ws.Range("D4").Value = "0dfee462-b107-462e-a961-d2072a9c7052"

$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMTgyODgzNH0.qnfaghQuQ8urPBlfTxcDVifN-5vWZry6DQG3hz7960k"
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMTgyODgzNn0.qodAJF06BHXed3DXfOcXs-VW4LwKOtHGj-CDM_vAb6c"
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDE4Mjg4Mzd9.vD6daaf1lKSwLPx863OJYyA98nlSfWZNM3Et_TiAXU4"
